$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Master")

# Row 94 - DfT Group row, add Rail / HSMRPG / AMIS example values
$ws.Range("B94").Value = "Rail"
$ws.Range("C94").Value = "Rail"
$ws.Range("D94").Value = "HSMRPG"
$ws.Range("E94").Value = "AMIS"
$ws.Range("F94").Value = "Rail"

# Row 95 - IPDC approval point row, add FBC / OBC / SOBC / pre-SOBC example values
$ws.Range("A95").Value = "IPDC approval point"
$ws.Range("B95").Value = "FBC"
$ws.Range("C95").Value = "OBC"
$ws.Range("D95").Value = "SOBC"
$ws.Range("E95").Value = "pre-SOBC"
$ws.Range("F95").Value = "FBC"

$ws.Range("F95").Select()
